$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws 'D2' '62.813.35'
Set-TextValue $ws 'E2' '  -1.65%  '

Set-TextValue $ws 'D3' '2.678.53'
Set-TextValue $ws 'E3' '  -2.14%  '

Set-TextValue $ws 'E4' '  +0.09%  '

Set-TextValue $ws 'D5' '553.99'
Set-TextValue $ws 'E5' '  -1.85%  '

Set-TextValue $ws 'D6' '158.20'
Set-TextValue $ws 'E6' '  -0.70%  '

Set-TextValue $ws 'D7' '1.00'
Set-TextValue $ws 'E7' '  +0.08%  '

Set-TextValue $ws 'D8' '0.592'
Set-TextValue $ws 'E8' '  -0.60%  '

Set-TextValue $ws 'E9' '  -2.83%  '

Set-TextValue $ws 'E10' '  -2.51%  '

Set-TextValue $ws 'D11' '0.369'
Set-TextValue $ws 'E11' '  -3.11%  '

Set-TextValue $ws 'D12' '5.38'
Set-TextValue $ws 'E12' '  -4.77%  '

Set-TextValue $ws 'D13' '3.153.07'
Set-TextValue $ws 'E13' '  -2.12%  '

Set-TextValue $ws 'D14' '26.46'
Set-TextValue $ws 'E14' '  -1.42%  '

Set-TextValue $ws 'D15' '62.750.74'
Set-TextValue $ws 'E15' '  -1.50%  '

Set-TextValue $ws 'D16' '0.0000147'
Set-TextValue $ws 'E16' '  -1.60%  '

Set-TextValue $ws 'D17' '2.685.73'
Set-TextValue $ws 'E17' '  -1.92%  '

Set-TextValue $ws 'D18' '11.90'
Set-TextValue $ws 'E18' '  -4.03%  '

Set-TextValue $ws 'D19' '4.63'
Set-TextValue $ws 'E19' '  -3.08%  '

Set-TextValue $ws 'D20' '344.51'
Set-TextValue $ws 'E20' '  -2.63%  '

Set-TextValue $ws 'D21' '6.29'
Set-TextValue $ws 'E21' '  -4.61%  '

Set-TextValue $ws 'E22' '  -0.15%  '

Set-TextValue $ws 'D23' '0.508'
Set-TextValue $ws 'E23' '  -2.45%  '

Set-TextValue $ws 'D24' '63.40'
Set-TextValue $ws 'E24' '  -1.20%  '

Set-TextValue $ws 'D25' '0.168'
Set-TextValue $ws 'E25' '  -0.58%  '

Set-TextValue $ws 'D26' '0.999'
Set-TextValue $ws 'E26' '  -0.23%  '

Set-TextValue $ws 'D27' '8.25'
Set-TextValue $ws 'E27' '  -1.22%  '

Set-TextValue $ws 'D28' '1.41'
Set-TextValue $ws 'E28' '  +8.06%  '

Set-TextValue $ws 'D29' '0.0₃0853'
Set-TextValue $ws 'E29' '  -5.24%  '

Set-TextValue $ws 'D30' '7.23'
Set-TextValue $ws 'E30' '  +0.67%  '

Set-TextValue $ws 'D31' '1.94'
Set-TextValue $ws 'E31' '  -0.84%  '

Set-TextValue $ws 'D32' '163.45'
Set-TextValue $ws 'E32' '  -0.14%  '

$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D33' '4.92'
Set-TextValue $ws 'E33' '  +0.82%  '

$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws 'D34' '1.49'
Set-TextValue $ws 'E34' '  +0.98%  '

Set-TextValue $ws 'E35' '  -0.02%  '

Set-TextValue $ws 'D36' '19.49'
Set-TextValue $ws 'E36' '  -2.91%  '

Set-TextValue $ws 'D37' '1.78'
Set-TextValue $ws 'E37' '  -1.08%  '

Set-TextValue $ws 'D38' '352.58'
Set-TextValue $ws 'E38' '  +2.54%  '

$ws.Range('B39').Value = 'SuiNetwork'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws 'D39' '0.949'
Set-TextValue $ws 'E39' '  -3.24%  '

$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D40' '6.23'
Set-TextValue $ws 'E40' '  -0.77%  '

Set-TextValue $ws 'D41' '3.99'
Set-TextValue $ws 'E41' '  -2.02%  '

Set-TextValue $ws 'D42' '38.48'
Set-TextValue $ws 'E42' '  +0.00%  '

Set-TextValue $ws 'D43' '20.91'
Set-TextValue $ws 'E43' '  -3.75%  '

Set-TextValue $ws 'D44' '20.20'
Set-TextValue $ws 'E44' '  -3.65%  '

Set-TextValue $ws 'D45' '0.617'
Set-TextValue $ws 'E45' '  -0.65%  '

Set-TextValue $ws 'D46' '0.0561'
Set-TextValue $ws 'E46' '  -3.29%  '

Set-TextValue $ws 'D47' '1.00'
Set-TextValue $ws 'E47' '  +0.12%  '

Set-TextValue $ws 'D48' '11.02'
Set-TextValue $ws 'E48' '  -0.33%  '

$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws 'D49' '0.0973'
Set-TextValue $ws 'E49' '  -2.69%  '

$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D50' '0.0243'
Set-TextValue $ws 'E50' '  -2.89%  '

Set-TextValue $ws 'D51' '128.52'
Set-TextValue $ws 'E51' '  -4.71%  '
